# Trade #69 closed at 2026-02-17 21:12:14 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: refreshed aggregate metrics after the new trade.
#  - Strategy Status sheet: MarketMaking row trade count / win rate refresh.
#  - All Trades sheet: trade #97 (row 98) closes out; new trade #130 appended (row 131).
#  - MarketMaking sheet: mirrors the same two trade updates (row 65 / new row 98).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal text value into a cell without Excel re-typing a
# date-shaped string ("2026-02-17") into a date serial number. Apply a Text
# number format before the write, then snap the style back to Normal so the
# on-disk cell ends up with no explicit style index (matching how the rest
# of the sheet stores its plain-text date column).
# ---------------------------------------------------------------------------
function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Materialise a present-but-empty cell (the source rows store these as a
# self-closed <c r="..." t="inlineStr"/> rather than omitting the cell
# entirely, e.g. an OPEN trade's not-yet-known Exit Price / Exit Reason).
function Set-BlankCell {
    param($range)
    $range.NumberFormat = "@"
    $range.Value = ""
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.18      # Total P&L %
$summary.Range("B6").Value = 97        # Total Trades
$summary.Range("B9").Value = 47.42     # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (row 5 = MarketMaking)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D5").Value = 64         # Trades
$status.Range("G5").Value = 50         # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Trade #97 (row 98) transitions from OPEN to CLOSED via early_exit.
$allTrades.Range("G98").Value = 0.15021
$allTrades.Range("H98").Value = "CLOSED"
$allTrades.Range("I98").Value = 0.1402
$allTrades.Range("K98").Value = 101.1
$allTrades.Range("L98").Value = "early_exit"
$allTrades.Range("M98").Value = 0.12

# New trade #130 (row 131), still OPEN.
$allTrades.Range("A131").Value = 130
Set-TextValue $allTrades.Range("B131") "2026-02-17"
$allTrades.Range("C131").Value = "21:12:08"
$allTrades.Range("D131").Value = "MarketMaking"
$allTrades.Range("E131").Value = "DOWN"
$allTrades.Range("F131").Value = 0.15
Set-BlankCell $allTrades.Range("G131")
$allTrades.Range("H131").Value = "OPEN"
$allTrades.Range("I131").Value = 0
$allTrades.Range("J131").Value = 0
$allTrades.Range("K131").Value = 101.0994048109029
Set-BlankCell $allTrades.Range("L131")
$allTrades.Range("M131").Value = 0
$allTrades.Range("N131").Value = 0
$allTrades.Range("O131").Value = 0
$allTrades.Range("P131").Value = 0.6
$allTrades.Range("Q131").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# MarketMaking sheet (same two trades, different column layout)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")

# Trade #97 (row 65) transitions from OPEN to CLOSED via early_exit.
$mm.Range("G65").Value = 0.15021
$mm.Range("H65").Value = "CLOSED"
$mm.Range("I65").Value = 0.1402
$mm.Range("K65").Value = 101.1
$mm.Range("P65").Value = "early_exit"
$mm.Range("Q65").Value = 0.12

# New trade #130 (row 98), still OPEN.
$mm.Range("A98").Value = 130
Set-TextValue $mm.Range("B98") "2026-02-17"
$mm.Range("C98").Value = "21:12:08"
$mm.Range("D98").Value = "MarketMaking"
$mm.Range("E98").Value = "DOWN"
$mm.Range("F98").Value = 0.15
Set-BlankCell $mm.Range("G98")
$mm.Range("H98").Value = "OPEN"
$mm.Range("I98").Value = 0
$mm.Range("J98").Value = 0
$mm.Range("K98").Value = 101.0994048109029
$mm.Range("L98").Value = 0
$mm.Range("M98").Value = 0
$mm.Range("N98").Value = 0.6
$mm.Range("O98").Value = "Normal spread capture: 19600 bps"
Set-BlankCell $mm.Range("P98")
$mm.Range("Q98").Value = 0
